# ZDD Bb4A Fundamentalstestcases.xlsx - apply "Add files via upload" edit
# Target sheet is "-the-name-of-your- module" (2nd worksheet, the ActiveSheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("-the-name-of-your- module")

# ---------------------------------------------------------------------
# 1) G18 header cell: tester name / date placeholder -> actual entry
# ---------------------------------------------------------------------
$ws.Range("G18").Value = "Drashti Desai `n11 august "

# ---------------------------------------------------------------------
# 2) Row heights that changed for this section of the sheet
# ---------------------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 73.2
$ws.Rows.Item(18).RowHeight = 55.8

# ---------------------------------------------------------------------
# 3) Fill in the previously-empty rows 20-24 with new test case data.
#    Columns A/B/F reuse existing text already used elsewhere on the
#    sheet; C/D hold the new sample run transcripts.
# ---------------------------------------------------------------------

# Row 24 ("Exiting with q")
$ws.Range("A24").Value = "Exiting with q"
$ws.Range("B24").Value = "positive "
$ws.Range("D24").Value = "*** Start of Measuring Strings Demo ***`nType a string (q - to quit):`nq`n*** End of Measuring Strings Demo ***"
$ws.Range("C24").Value = """q"""
$ws.Range("F24").Value = "PASS"

# Row 23 ("maximal edge case")
$ws.Range("A23").Value = "maximal edge case"
$ws.Range("B23").Value = "positive "
$ws.Range("D23").Value = "`nType a string (q - to quit):`nThis week I have a exams`nThe length of 'This week I have exams' is 25 characters"
$ws.Range("C23").Value = """This week I have a exams"""
$ws.Range("F23").Value = "PASS"

# Row 21 (" typical case")
$ws.Range("A21").Value = " typical case"
$ws.Range("B21").Value = "positive "
$ws.Range("D21").Value = "Type a string (q - to quit):`nHell is jail`nThe length of 'Hello' is 12 characters"
$ws.Range("C21").Value = """Hell is jail"""
$ws.Range("F21").Value = "PASS"

# Row 20 ("Handle exit input")
$ws.Range("A20").Value = "Handle exit input"
$ws.Range("B20").Value = "positive "
$ws.Range("D20").Value = "`nType a string (q - to quit):`nhey`nThe length of 'hey' is 3 characters`nType a string (q - to quit):`nq"
$ws.Range("C20").Value = """hey"",""q"""
$ws.Range("F20").Value = "PASS"

# Row 22 ("minimal edge case")
$ws.Range("A22").Value = "minimal edge case"
$ws.Range("B22").Value = "positive "
$ws.Range("D22").Value = "`n*** Start of Measuring Strings Demo ***`nType a string (q - to quit):`nThe length of '' is 0 characters"
$ws.Range("C22").Value = """"""
$ws.Range("F22").Value = "PASS"

# Row heights for the new data rows
$ws.Rows.Item(20).RowHeight = 92.4
$ws.Rows.Item(21).RowHeight = 52.8
$ws.Rows.Item(22).RowHeight = 79.2
$ws.Rows.Item(23).RowHeight = 66
$ws.Rows.Item(24).RowHeight = 79.2

# ---------------------------------------------------------------------
# 4) Update the sheet view: selection moves to E22, scroll position
#    resets (no pinned top-left cell)
# ---------------------------------------------------------------------
$ws.Range("E22").Select()

$wb.Save()
